$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match formatting of the existing header row (bold, border, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in team record values (Wins=77, Losses=85, Ties=0) for every data row
$ws.Range("AD2:AD66").Value = 77
$ws.Range("AE2:AE66").Value = 85
$ws.Range("AF2:AF66").Value = 0

$excel.CutCopyMode = 0
